$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "initial"
$ws2 = $wb.Worksheets.Item(2)   # "line_imp"

# ---------------------------------------------------------------------------
# Sheet "initial": one Newton-Raphson iteration worth of bus data updates
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = 0.5

$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 0.5
$ws1.Range("F3").Value = 0.5

$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 0
$ws1.Range("E4").ClearContents()
$ws1.Range("F4").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "line_imp": updated line impedance values
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = 0.05
$ws2.Range("C2").Value = 0.2

$ws2.Range("C3").Value = 0.1

$ws2.Range("B4").Value = 0.05
$ws2.Range("C4").Value = 0.15

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the saved UI state:
# "initial" keeps a stale selection at F4 (no longer the tab in view), while
# "line_imp" becomes the active, selected tab with its cursor at C5.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F4").Select()

$ws2.Activate()
$ws2.Range("C5").Select()

$win = $wb.Windows.Item(1)
$win.Left = 4680
